$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Date"
$ws.Range("A5").Value = "Author"

$ws.Range("A6").Select()
